$wb = $excel.ActiveWorkbook

# Update Table_Names sheet: A3 "TestTable3" -> "T4"
$wsTables = $wb.Worksheets.Item("Table_Names")
$wsTables.Range("A3").Value = "T4"

# Update Field_Names sheet: clear out TestValue3.1..3.4 in A6:A9
$wsFields = $wb.Worksheets.Item("Field_Names")
$wsFields.Range("A6:A9").ClearContents()
